# Update "想去人数" (F column) figures to the latest scraped counts.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 5901
$ws.Range("F9").Value = 1376
$ws.Range("F12").Value = 1962
$ws.Range("F17").Value = 179
$ws.Range("F22").Value = 66
$ws.Range("F23").Value = 3684
$ws.Range("F25").Value = 2925
$ws.Range("F27").Value = 2478
$ws.Range("F28").Value = 4245
$ws.Range("F32").Value = 1343
$ws.Range("F33").Value = 99
$ws.Range("F36").Value = 30

# Sheet "演出" (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 918
$ws.Range("F17").Value = 13

# Sheet "全部类型" (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5901
$ws.Range("F8").Value = 1376
$ws.Range("F11").Value = 1962
$ws.Range("F15").Value = 918
$ws.Range("F18").Value = 179
$ws.Range("F21").Value = 3684
$ws.Range("F22").Value = 13
$ws.Range("F25").Value = 2925
$ws.Range("F26").Value = 2478
$ws.Range("F27").Value = 4245
$ws.Range("F30").Value = 1343
